$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 174.55
$ws.Range("I2").Value = 183.94444
$ws.Range("K2").Value = 183.94444
$ws.Range("M2").Value = -70.94443999999999
$ws.Range("H29").Value = 699.5
$ws.Range("J29").Value = 699.5
$ws.Range("L29").Value = 2098.5
$ws.Range("N29").Value = -2660.5
$ws.Range("H32").Value = 4876.6665
$ws.Range("I32").Value = 6983.3335
$ws.Range("J32").Value = 4350
$ws.Range("K32").Value = 6983.3335
$ws.Range("L32").Value = 4350
$ws.Range("M32").Value = -6657.3335
$ws.Range("N32").Value = -5002
$ws.Range("H38").Value = 106.23077
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("H39").Value = 605.25
$ws.Range("I39").Value = 158.5
$ws.Range("K39").Value = 475.5
$ws.Range("M39").Value = -179.5
$ws.Range("H43").Value = 2486.25
$ws.Range("I43").Value = 2570
$ws.Range("K43").Value = 2570
$ws.Range("M43").Value = -2501
$ws.Range("H49").Value = 50
$ws.Range("I49").Value = 50
$ws.Range("K49").Value = 150
$ws.Range("H98").Value = 4634.511
$ws.Range("I98").Value = 4799.8975
$ws.Range("J98").Value = 3559.5
$ws.Range("K98").Value = 4799.8975
$ws.Range("L98").Value = 3559.5
$ws.Range("M98").Value = -3301.8975
$ws.Range("N98").Value = -6555.5
$ws.Range("H115").Value = 300
$ws.Range("I115").Value = 300
$ws.Range("K115").Value = 900
$ws.Range("M115").Value = 667
$ws.Range("H122").Value = 4634.511
$ws.Range("I122").Value = 4799.8975
$ws.Range("J122").Value = 3559.5
$ws.Range("K122").Value = 14399.6925
$ws.Range("L122").Value = 10678.5
$ws.Range("M122").Value = -11949.6925
$ws.Range("N122").Value = -15578.5
$ws.Range("H131").Value = 412500
$ws.Range("I131").Value = 672833.3
$ws.Range("K131").Value = 2018499.9
$ws.Range("M131").Value = -2013459.9
$ws.Range("H132").Value = 3770.8572
$ws.Range("I132").Value = 3849.2646
$ws.Range("J132").Value = 1105
$ws.Range("K132").Value = 11547.7938
$ws.Range("L132").Value = 3315
$ws.Range("M132").Value = -9017.793799999999
$ws.Range("N132").Value = -8375
$ws.Range("H138").Value = 1796.7
$ws.Range("I138").Value = 1255.8214
$ws.Range("J138").Value = 3058.75
$ws.Range("K138").Value = 3767.4642
$ws.Range("L138").Value = 9176.25
$ws.Range("M138").Value = 1372.5358
$ws.Range("N138").Value = -19456.25
$ws.Range("H141").Value = 6835.952
$ws.Range("I141").Value = 6427.8
$ws.Range("J141").Value = 14999
$ws.Range("K141").Value = 19283.4
$ws.Range("L141").Value = 44997
$ws.Range("M141").Value = -14103.4
$ws.Range("N141").Value = -55357
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4068.2144
$ws.Range("I45").Value = 2863.3462
$ws.Range("J45").Value = 6026.125
$ws.Range("K45").Value = 2863.3462
$ws.Range("L45").Value = 6026.125
$ws.Range("M45").Value = -2486.3462
$ws.Range("N45").Value = -6780.125
$ws.Range("H74").Value = 8666.5
$ws.Range("I74").Value = 7374.75
$ws.Range("K74").Value = 7374.75
$ws.Range("M74").Value = -6500.75
$ws.Range("H77").Value = 8666.5
$ws.Range("I77").Value = 7374.75
$ws.Range("K77").Value = 36873.75
$ws.Range("M77").Value = -32505.75
$ws.Range("H132").Value = 4010.4385
$ws.Range("I132").Value = 3228.8333
$ws.Range("J132").Value = 6198.933
$ws.Range("K132").Value = 9686.499899999999
$ws.Range("L132").Value = 18596.799
$ws.Range("M132").Value = -7156.499899999999
$ws.Range("N132").Value = -23656.799
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4046.4707
$ws.Range("I105").Value = 4099.4375
$ws.Range("K105").Value = 4099.4375
$ws.Range("M105").Value = -2352.4375
$ws.Range("H134").Value = 5397.7915
$ws.Range("I134").Value = 5483.9546
$ws.Range("J134").Value = 4450
$ws.Range("K134").Value = 16451.8638
$ws.Range("L134").Value = 13350
$ws.Range("M134").Value = -13916.8638
$ws.Range("N134").Value = -18420
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 5555
$ws.Range("I3").Value = 5555
$ws.Range("K3").Value = 5555
$ws.Range("M3").Value = -5442
$ws.Range("H7").Value = 138.5
$ws.Range("J7").Value = 221.66667
$ws.Range("L7").Value = 221.66667
$ws.Range("N7").Value = -447.66667
$ws.Range("H22").Value = 4276.7
$ws.Range("I22").Value = 5895
$ws.Range("J22").Value = 500.66666
$ws.Range("K22").Value = 5895
$ws.Range("L22").Value = 500.66666
$ws.Range("M22").Value = -5545
$ws.Range("N22").Value = -1200.66666
$ws.Range("H137").Value = 84435
$ws.Range("J137").Value = 84435
$ws.Range("L137").Value = 84435
$ws.Range("N137").Value = -94635
$ws.Range("H139").Value = 137496.42
$ws.Range("I139").Value = 129000
$ws.Range("K139").Value = 129000
$ws.Range("M139").Value = -123860
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I4").Value = 5618676.5
$ws.Range("K4").Value = 16856029.5
$ws.Range("M4").Value = -16855917.5
$ws.Range("H14").Value = 332.57144
$ws.Range("I14").Value = 332.57144
$ws.Range("K14").Value = 997.71432
$ws.Range("M14").Value = -824.71432
$ws.Range("H23").Value = 11234.889
$ws.Range("J23").Value = 33503
$ws.Range("L23").Value = 100509
$ws.Range("N23").Value = -100979
$ws.Range("H29").Value = 28388.375
$ws.Range("I29").Value = 322.83334
$ws.Range("J29").Value = 112585
$ws.Range("K29").Value = 968.5000200000001
$ws.Range("L29").Value = 337755
$ws.Range("M29").Value = -691.5000200000001
$ws.Range("N29").Value = -338309
$ws.Range("H33").Value = 693.61536
$ws.Range("J33").Value = 1358.1666
$ws.Range("L33").Value = 8148.9996
$ws.Range("N33").Value = -8714.999599999999
$ws.Range("H51").Value = 60499.5
$ws.Range("I51").Value = 999
$ws.Range("K51").Value = 2997
$ws.Range("M51").Value = -2537
$ws.Range("H96").Value = 5249.5
$ws.Range("J96").Value = 9999
$ws.Range("L96").Value = 29997
$ws.Range("N96").Value = -34115
$ws.Range("H97").Value = 1395.1305
$ws.Range("I97").Value = 1326.8
$ws.Range("J97").Value = 1523.25
$ws.Range("K97").Value = 3980.4
$ws.Range("L97").Value = 4569.75
$ws.Range("M97").Value = -3484.4
$ws.Range("N97").Value = -5561.75
$ws.Range("H98").Value = 347
$ws.Range("I98").Value = 350
$ws.Range("J98").Value = 344
$ws.Range("K98").Value = 1050
$ws.Range("L98").Value = 1032
$ws.Range("M98").Value = 448
$ws.Range("H104").Value = 4289
$ws.Range("J104").Value = 5833
$ws.Range("L104").Value = 17499
$ws.Range("N104").Value = -22741
$ws.Range("H112").Value = 3666.3333
$ws.Range("I112").Value = 2999.5
$ws.Range("K112").Value = 8998.5
$ws.Range("M112").Value = -7890.5
$ws.Range("H132").Value = 83334310
$ws.Range("J132").Value = 1248.25
$ws.Range("L132").Value = 11234.25
$ws.Range("N132").Value = -16294.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2965.923
$ws.Range("I80").Value = 2618.5
$ws.Range("J80").Value = 4124
$ws.Range("K80").Value = 2618.5
$ws.Range("L80").Value = 4124
$ws.Range("M80").Value = -1620.5
$ws.Range("N80").Value = -6120
$ws.Range("H83").Value = 2965.923
$ws.Range("I83").Value = 2618.5
$ws.Range("J83").Value = 4124
$ws.Range("K83").Value = 13092.5
$ws.Range("L83").Value = 20620
$ws.Range("M83").Value = -8100.5
$ws.Range("N83").Value = -30604
$ws.Range("H122").Value = 2170.875
$ws.Range("J122").Value = 1751.25
$ws.Range("L122").Value = 5253.75
$ws.Range("N122").Value = -10153.75
$ws.Range("H132").Value = 3025.9
$ws.Range("I132").Value = 2030.5294
$ws.Range("K132").Value = 6091.5882
$ws.Range("M132").Value = -3561.5882
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2525.4375
$ws.Range("I7").Value = 2460.4666
$ws.Range("K7").Value = 2460.4666
$ws.Range("M7").Value = -2348.4666
$ws.Range("H16").Value = 2944.9583
$ws.Range("I16").Value = 2651.3809
$ws.Range("K16").Value = 2651.3809
$ws.Range("M16").Value = -2481.3809
$ws.Range("H40").Value = 4426.8184
$ws.Range("I40").Value = 3528.4285
$ws.Range("K40").Value = 3528.4285
$ws.Range("M40").Value = -3392.4285
$ws.Range("H46").Value = 18780
$ws.Range("I46").Value = 2699.75
$ws.Range("J46").Value = 50940.5
$ws.Range("K46").Value = 2699.75
$ws.Range("L46").Value = 50940.5
$ws.Range("M46").Value = -2511.75
$ws.Range("N46").Value = -51316.5
$ws.Range("H126").Value = 2525.4375
$ws.Range("I126").Value = 2460.4666
$ws.Range("K126").Value = 7381.399800000001
$ws.Range("M126").Value = -4911.399800000001
$ws.Range("H132").Value = 20778.428
$ws.Range("I132").Value = 20778.428
$ws.Range("K132").Value = 62335.284
$ws.Range("M132").Value = -59805.284
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 602.7143
$ws.Range("I107").Value = 638.7857
$ws.Range("J107").Value = 530.5714
$ws.Range("K107").Value = 1916.3571
$ws.Range("L107").Value = 1591.7142
$ws.Range("M107").Value = 3.642899999999827
$ws.Range("N107").Value = -5431.7142
$ws.Range("H122").Value = 3225.3809
$ws.Range("I122").Value = 2422.5
$ws.Range("J122").Value = 5794.6
$ws.Range("K122").Value = 7267.5
$ws.Range("L122").Value = 17383.8
$ws.Range("M122").Value = -4817.5
$ws.Range("N122").Value = -22283.8
$ws.Range("H132").Value = 9490
$ws.Range("I132").Value = 7483.3335
$ws.Range("K132").Value = 22450.0005
$ws.Range("M132").Value = -19920.0005
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N38").ClearContents()
$ws.Range("M49").Value = -14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N98").Value = -4028
